$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.884.84'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').Value = '2.420.03'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '551.04'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.49'
$ws.Range('E6').Value = '  +0.75%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  +2.43%  '
$ws.Range('E9').Value = '  -2.29%  '
$ws.Range('E10').Value = '  -3.49%  '
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('E12').Value = '  -2.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '25.48'
$ws.Range('E13').Value = '  +3.31%  '
$ws.Range('D14').Value = '2.851.44'
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('D15').Value = '59.804.92'
$ws.Range('E15').Value = '  +0.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000137'
$ws.Range('E16').Value = '  -1.76%  '
$ws.Range('D17').Value = '2.394.68'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '330.18'
$ws.Range('E20').Value = '  -1.89%  '
$ws.Range('E21').Value = '  -4.52%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.64'
$ws.Range('E23').Value = '  +3.10%  '
$ws.Range('E24').Value = '  +1.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.74'
$ws.Range('E25').Value = '  +3.40%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.37'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').Value = '0.0₃0776'
$ws.Range('E28').Value = '  +1.34%  '
$ws.Range('E29').Value = '  -1.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '169.10'
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('E31').Value = '  -2.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.70'
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  +0.27%  '
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('E38').Value = '  -2.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '39.54'
$ws.Range('E39').Value = '  -1.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.412'
$ws.Range('E40').Value = '  -2.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '314.37'
$ws.Range('E41').Value = '  +6.14%  '
$ws.Range('E42').Value = '  -2.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '139.42'
$ws.Range('E43').Value = '  -2.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0967'
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.58'
$ws.Range('E46').Value = '  +1.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.578'
$ws.Range('E47').Value = '  +1.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0225'
$ws.Range('E48').Value = '  -0.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.389'
$ws.Range('E49').Value = '  -2.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.61'
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.05'
$ws.Range('E51').Value = '  +0.12%  '
